# The sheet contains a weekly price log for "Zanahoria" (carrot) at the
# "Terminal La Palmera de La Serena" market. A new weekly record needs to be
# inserted at the top of the data block (row 94), pushing all the existing
# records (rows 94-244) down by one row (to rows 95-245), and the new row 94
# is populated with this week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 94; this shifts rows 94:244 down to 95:245
# and Excel automatically grows the sheet dimension to A1:R245.
$ws.Rows.Item(94).Insert()

# Populate the newly inserted row 94 with the new weekly record.
$ws.Cells.Item(94, 1).Value  = 8
$ws.Cells.Item(94, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(94, 3).Value  = "Coquimbo"
$ws.Cells.Item(94, 4).Value  = (Get-Date -Year 2021 -Month 12 -Day 10 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(94, 5).Value  = 4
$ws.Cells.Item(94, 6).Value  = 100114013
$ws.Cells.Item(94, 7).Value  = "Zanahoria"
$ws.Cells.Item(94, 8).Value  = "Sin especificar"
$ws.Cells.Item(94, 9).Value  = "Primera"
$ws.Cells.Item(94, 10).Value = 600
$ws.Cells.Item(94, 11).Value = 6000
$ws.Cells.Item(94, 12).Value = 7000
$ws.Cells.Item(94, 13).Value = 6500
$ws.Cells.Item(94, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(94, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(94, 16).Value = 325
$ws.Cells.Item(94, 17).Value = 20
$ws.Cells.Item(94, 18).Value = "Hortaliza"

# Apply the same date number format used by the rest of column D so the new
# row matches the formatting of all the other data rows.
$ws.Cells.Item(94, 4).NumberFormat = $ws.Cells.Item(95, 4).NumberFormat
